$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Before the edit: the last data row is 30 and row 31 is the "Yht" (total)
# row holding the label + =SUM(C6:C30) formula.
# After the edit: 4 new data rows are appended (31-34) and the total row
# moves down to row 35 with its formula extended to =SUM(C6:C34).

# --- Move the total row from 31 down to 35 -------------------------------
# Copy its formatting first (reusing the existing style records instead of
# creating new duplicate ones), then copy over its label/formula content.
$ws.Range("B31:D31").Copy()
$ws.Range("B35:D35").PasteSpecial(-4122)   # xlPasteFormats
$excel.CutCopyMode = $false

$ws.Cells.Item(35, 2).Value = "Yht"
$ws.Rows.Item(35).RowHeight = 18.75

# --- Update the existing last data row (C30: 4 -> 2) ----------------------
$ws.Cells.Item(30, 3).Value = 2

# --- Fill in the new data rows 31-34 --------------------------------------
$ws.Cells.Item(31, 2).Value = 45352
$ws.Cells.Item(31, 3).Value = 1
$ws.Cells.Item(31, 4).Value = "Koitin tehdä työkaluista enemmän dynaamisia. Onnistuin ideoimaan hyvän implementointi metodin."
$ws.Rows.Item(31).RowHeight = 37.5

$ws.Cells.Item(32, 2).Value = 45355
$ws.Cells.Item(32, 3).Value = 4
$ws.Cells.Item(32, 4).Value = "Muutin navigaatiopalkin luomisen enemmän moduulariksi. Vaatii vielä vähän hiomista asetusten ja joidenkin tiedostojen suhteen."
$ws.Rows.Item(32).RowHeight = 56.25

$ws.Cells.Item(33, 2).Value = 45356
$ws.Cells.Item(33, 3).Value = 3
$ws.Cells.Item(33, 4).Value = "Koitin saada modulaarisen navigaatiopalkin toimimaan täysin oikein. En saanut yhtä tiedostoa toimimaan oikein."
$ws.Rows.Item(33).RowHeight = 56.25

$ws.Cells.Item(34, 2).Value = 45357
$ws.Cells.Item(34, 3).Value = 5
$ws.Cells.Item(34, 4).Value = "Sain modulaarisen navigaation toimimaan kunolla. Korjasin gitignore"
$ws.Rows.Item(34).RowHeight = 37.5

# Apply the same visual style as the other data rows (copy formats from
# row 30, the last original data row) onto the 4 newly filled rows so they
# reuse the existing B/C/D style records rather than creating new ones.
$ws.Range("B30:D30").Copy()
$ws.Range("B31:D34").PasteSpecial(-4122)   # xlPasteFormats
$excel.CutCopyMode = $false

# --- Fix up the total formula (now on row 35) to cover the new range -----
$ws.Cells.Item(35, 3).Formula = "=SUM(C6:C34)"

# --- Match the saved sheet view state (selection) -------------------------
$ws.Range("D40").Select()
